$p = $ppt.ActivePresentation

# --- Slide 2: title "Summary of CFA Issues & Resolution (1)" -> "Summary of Issues & Resolution (1)" ---
$s2 = $p.Slides.Item(2)
$titleTr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$titleRun2 = $titleTr2.Characters(8, 17)
$titleRun2.Text = " of Issues & "

# --- Slide 2: "...to help deciding if similar..." -> "...to help decide if similar..." ---
$bodyTr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$outcomePara = $bodyTr2.Paragraphs(4, 1)
$outcomeLen = $outcomePara.Text.Length
$outcomeRun = $outcomePara.Characters(1, $outcomeLen)
$outcomeRun.Text = "The outcome is to keep the text in the realization I-D + Add NEW Scope text to both I-Ds to help decide if similar issues are raised in the future"

# --- Slide 5: title "Summary of CFA Issues & Resolution (2)" -> "Summary of Issues & Resolution (2)" ---
$s5 = $p.Slides.Item(5)
$titleTr5 = $s5.Shapes.Item(1).TextFrame.TextRange
$titleRun5 = $titleTr5.Characters(8, 17)
$titleRun5.Text = " of Issues & "

# --- Slide 7: "The authors think that content is almost stable"
#     -> 3 runs: "The authors think " / "that the content " / "is almost stable" ---
$s7 = $p.Slides.Item(7)
$bodyTr7 = $s7.Shapes.Item(2).TextFrame.TextRange
$stablePara = $bodyTr7.Paragraphs(1, 1)
$midRun = $stablePara.Characters(19, 13)
$midRun.Text = "that the content "
